# Adds bold section headings (Introduction / Paragraph 1-3 / Conclusion) to the
# essay template and re-styles the {question} paragraph as a large, bold,
# centered title. Also drops the trailing space that used to follow each
# placeholder token now that it sits directly under its own heading.

$d = $word.ActiveDocument

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$xmlFooter = '</w:body></w:document>'

function Set-PlainBody($index, $text) {
    # Re-writes a placeholder paragraph (by 1-based index) so its run text
    # has no trailing space, keeping the original (non-bold) run formatting.
    $paragraph = $d.Paragraphs.Item($index)
    $xml = $xmlHeader + '<w:p><w:pPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>' + $xmlFooter
    $paragraph.Range.InsertXML($xml)
}

function Insert-BoldHeading($index, $text) {
    # Inserts a new bold heading paragraph at (1-based) $index, pushing the
    # paragraph that used to live there down to $index + 1.
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphBefore()
    $headingPara = $d.Paragraphs.Item($index)
    $xml = $xmlHeader + '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>' + $xmlFooter
    $headingPara.Range.InsertXML($xml)
}

function Insert-IntroductionHeading($index) {
    # "Introduction:" is authored as two separate runs ("Introduction" + ":").
    $anchor = $d.Paragraphs.Item($index)
    $anchor.Range.InsertParagraphBefore()
    $headingPara = $d.Paragraphs.Item($index)
    $xml = $xmlHeader + '<w:p><w:pPr><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Introduction</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>:</w:t></w:r></w:p>' + $xmlFooter
    $headingPara.Range.InsertXML($xml)
}

# Work from the bottom of the document upward so the paragraph indices for
# the earlier placeholders stay valid while later ones are edited.

# --- Conclusion heading, ahead of the existing {conclusion} paragraph -----
Insert-BoldHeading 11 "Conclusion:"

# --- Paragraph 3 heading + {body_3} trailing-space cleanup ----------------
Set-PlainBody 9 "{body_3}"
Insert-BoldHeading 9 "Paragraph 3:"

# --- Paragraph 2 heading + {body_2} trailing-space cleanup ----------------
Set-PlainBody 7 "{body_2}"
Insert-BoldHeading 7 "Paragraph 2:"

# --- Paragraph 1 heading + {body_1} trailing-space cleanup ----------------
Set-PlainBody 5 "{body_1}"
Insert-BoldHeading 5 "Paragraph 1:"

# --- Introduction heading + {introduction} trailing-space cleanup ---------
Set-PlainBody 3 "{introduction}"
Insert-IntroductionHeading 3

# --- {question} becomes a large, bold, centered title ---------------------
$question = $d.Paragraphs.Item(1)
$questionXml = $xmlHeader + '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>{question}</w:t></w:r></w:p>' + $xmlFooter
$question.Range.InsertXML($questionXml)
